$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 4 (Famicom (Large) and Famicom Disk entries).
# This shifts the existing rows 4-1004 down to 6-1006, carrying their styles along,
# and grows the used range / blank filler rows automatically.
$ws.Rows("4:5").Insert()

# Row 4: Nintendo / Famicom (Large) / 91 x 67.5 / 13.5 / (blank) / (blank)
$ws.Cells.Item(4, 1).Value = "Nintendo"
$ws.Cells.Item(4, 2).Value = "Famicom (Large)"
$ws.Cells.Item(4, 3).Value = "91 x 67.5"
$ws.Cells.Item(4, 4).Value = 13.5

# Row 5: Nintendo / Famicom Disk / 52 x 22 / NA / (blank) / (blank)
$ws.Cells.Item(5, 1).Value = "Nintendo"
$ws.Cells.Item(5, 2).Value = "Famicom Disk"
$ws.Cells.Item(5, 3).Value = "52 x 22"
$ws.Cells.Item(5, 4).Value = "NA"

# Row 6 (was row 4, Gameboy): now confirmed -> F gets the "x" mark.
# Copy the "Confirmed" cell format (from a row that already has the x-style, e.g. F7)
# so the cell reuses the existing style instead of Excel fabricating a new one.
$ws.Cells.Item(7, 6).Copy() | Out-Null
$ws.Cells.Item(6, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6, 6).Value = "x"

# Row 8 (was row 6, N64): now confirmed -> F gets the "x" mark
$ws.Cells.Item(7, 6).Copy() | Out-Null
$ws.Cells.Item(8, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 6).Value = "x"

$excel.CutCopyMode = 0

# Row 10 (was row 8, NES (Back)): Top Fold is now "NA" instead of blank
$ws.Cells.Item(10, 4).Value = "NA"

# Row 12 (was row 10, SNES (Back)): Top Fold is now "NA" instead of blank
$ws.Cells.Item(12, 4).Value = "NA"
